$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Level (column C) updates ---
$ws.Range("C5").Value = 2
$ws.Range("C8").Value = 31
$ws.Range("C9").Value = 31
$ws.Range("C10").Value = 3
$ws.Range("C11").Value = 31
$ws.Range("C13").Value = 31
$ws.Range("C14").Value = 31
$ws.Range("C15").Value = 31
$ws.Range("C16").Value = 4
$ws.Range("C17").Value = 31
$ws.Range("C18").Value = 31

# --- New quest text for row 16 (Quest column D) ---
$ws.Range("D16").Value = "trees;3|grave;1"

# --- Selection matches the authored workbook view ---
[void]$ws.Range("D17").Select()
